# Fruta / hortaliza, semanal
#
# The underlying data rows (2..16) got reshuffled: for each destination
# row, the values in columns D (Fecha), L (Calidad), M (Volumen),
# N (Precio minimo), O (Precio maximo), P (Precio promedio ponderado),
# Q (Unidad de comercializacion) and S (Precio $/Kg) are replaced with
# the values that, prior to the edit, lived in a different row of the
# same block. All other columns (A, B, C, E-K, R, T) are identical
# across every row and therefore stay untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the "before" state of the columns that move, keyed by row number.
$cols = @("D", "L", "M", "N", "O", "P", "Q", "S")

$before = @{}
for ($r = 2; $r -le 16; $r++) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Range("$c$r").Value2()
    }
    $before[$r] = $rowVals
}

# Destination row -> source row (source row's data moves into destination row).
$mapping = @{
    2  = 14
    3  = 15
    4  = 9
    5  = 16
    6  = 7
    7  = 8
    8  = 3
    9  = 10
    10 = 5
    11 = 2
    12 = 4
    13 = 6
    14 = 11
    15 = 12
    16 = 13
}

foreach ($destRow in $mapping.Keys) {
    $srcRow = $mapping[$destRow]
    $srcVals = $before[$srcRow]
    foreach ($c in $cols) {
        $ws.Range("$c$destRow").Value2 = $srcVals[$c]
    }
}

$wb.Save()
